# Hortaliza, Vega Modelo de Temuco - Zanahoria
# Insert a new weekly price-report row at row 275, shifting the existing
# rows 275-334 down to 276-335 (dimension grows from A1:R334 to A1:R335).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 275 (pushes old 275..334 down to 276..335)
$ws.Rows.Item(275).Insert()

# Populate the newly inserted row 275 with the new data point
$ws.Range("A275").Value2 = 10
$ws.Range("B275").Value2 = "Vega Modelo de Temuco"
$ws.Range("C275").Value2 = "La Araucanía"
$ws.Range("D275").Value2 = 44782
$ws.Range("E275").Value2 = 9
$ws.Range("F275").Value2 = 100114013
$ws.Range("G275").Value2 = "Zanahoria"
$ws.Range("H275").Value2 = "Sin especificar"
$ws.Range("I275").Value2 = "Primera"
$ws.Range("J275").Value2 = 50
$ws.Range("K275").Value2 = 10000
$ws.Range("L275").Value2 = 10000
$ws.Range("M275").Value2 = 10000
$ws.Range("N275").Value2 = "$/saco 20 kilos"
$ws.Range("O275").Value2 = "Región de La Araucanía"
$ws.Range("P275").Value2 = 400
$ws.Range("Q275").Value2 = 25
$ws.Range("R275").Value2 = "Hortaliza"
